$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header (H1) onto the two new
# header cells so they pick up the same bold/border/center-top style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-27
$iValues = @(5,7,8,9,6,5,3,7,9,8,9,8,7,7,9,9,9,7,9,7,8,7,5,4,6,5)
$jValues = @(6,7,8,9,6,5,3,8,9,8,9,8,8,7,9,9,9,8,9,7,8,7,5,4,6,5)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
